$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "You Died" screen localization row ---
$ws.Range("A8").Value2 = "You Died"
$ws.Range("B8").Value2 = "You Died"
$ws.Range("C8").Value2 = "阵亡"

# --- "Play Again" / "Return Home" rows ---
$ws.Range("A9").Value2 = "Play Again"
$ws.Range("B9").Value2 = "Play Again"

$ws.Range("A10").Value2 = "Return Home"
$ws.Range("B10").Value2 = "Return Home"

$ws.Range("C10").Value2 = "返回主页"
$ws.Range("C9").Value2 = "再次游玩"

# --- Rename "Start" -> "Start Adventure" on row 2 ---
$ws.Range("C2").Value2 = "开始冒险"
$ws.Range("A2").Value2 = "Start Adventure"
$ws.Range("B2").Value2 = "Start Adventure"

# Column A (English key column) carries the same "Noto Sans" cell style as
# column B, for every populated data row. Copy the format from B2 (which
# already uses that style) instead of touching Font.Name directly, so we
# reuse the existing style entry instead of minting a new, unused one.
$ws.Range("B2").Copy()
$ws.Range("A2:A10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to mirror the new bottom-right of the table.
$ws.Range("D10").Select()
